$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Update result column: C3 stays "Y"; all others in C2:C7 become "N"
$ws.Range("C2").Value = "N"
$ws.Range("C4").Value = "N"
$ws.Range("C5").Value = "N"
$ws.Range("C6").Value = "N"
$ws.Range("C7").Value = "N"

# Update the active selection to reflect the changed range C4:C7
$ws.Range("C4:C7").Select()
